# Apply updates described by the diff: header rename + updated B/D/E values (rows 2-19)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename: E1 "strength (raw)" -> "strength (RMS)"
$ws.Range("E1").Value = "strength (RMS)"

# New values per row: reactionTime (B), difference (D), strength (E)
$data = @(
    @{ Row = 2;  B = 14;    D = 13.6;               E = 76.2 }
    @{ Row = 3;  B = 14.5;  D = 16.5;               E = 79.5 }
    @{ Row = 4;  B = 14;    D = 16.8;               E = 75.2 }
    @{ Row = 5;  B = 14;    D = 21.2;               E = 77.8 }
    @{ Row = 6;  B = 14;    D = 16.4;               E = 73.2 }
    @{ Row = 7;  B = 14.4;  D = 14;                 E = 67.2 }
    @{ Row = 8;  B = 14;    D = 16.4;               E = 72.59999999999999 }
    @{ Row = 9;  B = 14;    D = 18.8;               E = 73 }
    @{ Row = 10; B = 12.8;  D = 19.6;               E = 77.40000000000001 }
    @{ Row = 11; B = 14;    D = 31.6;               E = 71.8 }
    @{ Row = 12; B = 14.4;  D = 16.8;               E = 66.8 }
    @{ Row = 13; B = 14;    D = 17.6;               E = 71.40000000000001 }
    @{ Row = 14; B = 13.33; D = 17.33;              E = 73.67 }
    @{ Row = 15; B = 14;    D = 17;                 E = 66.25 }
    @{ Row = 16; B = 14;    D = 20.4;               E = 78.8 }
    @{ Row = 17; B = 13.6;  D = 20.8;               E = 82.2 }
    @{ Row = 18; B = 15.2;  D = 16;                 E = 69.59999999999999 }
    @{ Row = 19; B = 12.8;  D = 18;                 E = 72.8 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
}
